$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; existing rows 20..61 shift down to 21..62
$ws.Rows.Item(20).Insert()

# Fill in the new row 20 with data (same as the row that was previously at row 20,
# except for the changed fields described in the diff)
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C20").Value = "Los Lagos"
$ws.Range("D20").Value = (Get-Date -Year 2023 -Month 2 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100101
$ws.Range("H20").Value = "Berries"
$ws.Range("I20").Value = 100101001
$ws.Range("J20").Value = "Arándano (blue)"
$ws.Range("K20").Value = "Sin especificar"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 2000
$ws.Range("O20").Value = 2200
$ws.Range("P20").Value = 2100
$ws.Range("Q20").Value = "$/bandeja 2 kilos"
$ws.Range("R20").Value = "Provincia de Curicó"
$ws.Range("S20").Value = 1050
$ws.Range("T20").Value = 2
